# Atualizacao de bases das ligas, do dia: 03-03-2024 as 00:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 140
$ws.Cells.Item(140, 2).Value = 7493310
$ws.Cells.Item(140, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(140, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(140, 5).Value = 45259.8125
$ws.Cells.Item(140, 6).Value = "Libertad Asuncion"
$ws.Cells.Item(140, 7).Value = "Tacuary"
$ws.Cells.Item(140, 8).Value = 1
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = "A"
$ws.Cells.Item(140, 11).Value = 1.363
$ws.Cells.Item(140, 12).Value = 5
$ws.Cells.Item(140, 13).Value = 7
$ws.Cells.Item(140, 14).Value = 1.571
$ws.Cells.Item(140, 15).Value = 4.2
$ws.Cells.Item(140, 16).Value = 4.75
$ws.Cells.Item(140, 17).Value = -0.75
$ws.Cells.Item(140, 18).Value = 1.8
$ws.Cells.Item(140, 19).Value = 2
$ws.Cells.Item(140, 20).Value = 2.75
$ws.Cells.Item(140, 21).Value = 1.8
$ws.Cells.Item(140, 22).Value = 2
$ws.Cells.Item(140, 23).Value = -1
$ws.Cells.Item(140, 24).Value = -1
$ws.Cells.Item(140, 25).Value = 3.75
$ws.Cells.Item(140, 26).Value = -1
$ws.Cells.Item(140, 27).Value = 1
$ws.Cells.Item(140, 28).Value = 0.4
$ws.Cells.Item(140, 29).Value = -0.5

# Row 141
$ws.Cells.Item(141, 2).Value = 7493431
$ws.Cells.Item(141, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(141, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(141, 5).Value = 45259.8125
$ws.Cells.Item(141, 6).Value = "Sportivo Trinidense"
$ws.Cells.Item(141, 7).Value = "Guairena FC"
$ws.Cells.Item(141, 8).Value = 7
$ws.Cells.Item(141, 9).Value = 2
$ws.Cells.Item(141, 10).Value = "H"
$ws.Cells.Item(141, 11).Value = 2.05
$ws.Cells.Item(141, 12).Value = 3.3
$ws.Cells.Item(141, 13).Value = 3.3
$ws.Cells.Item(141, 14).Value = 2.6
$ws.Cells.Item(141, 15).Value = 3.1
$ws.Cells.Item(141, 16).Value = 2.6
$ws.Cells.Item(141, 17).Value = 0
$ws.Cells.Item(141, 18).Value = 1.925
$ws.Cells.Item(141, 19).Value = 1.875
$ws.Cells.Item(141, 20).Value = 2.5
$ws.Cells.Item(141, 21).Value = 2
$ws.Cells.Item(141, 22).Value = 1.8
$ws.Cells.Item(141, 23).Value = 1.6
$ws.Cells.Item(141, 24).Value = -1
$ws.Cells.Item(141, 25).Value = -1
$ws.Cells.Item(141, 26).Value = 0.925
$ws.Cells.Item(141, 27).Value = -1
$ws.Cells.Item(141, 28).Value = 1
$ws.Cells.Item(141, 29).Value = -1

# Row 143
$ws.Cells.Item(143, 2).Value = 7493433
$ws.Cells.Item(143, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(143, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(143, 5).Value = 45261.8125
$ws.Cells.Item(143, 6).Value = "Sportivo Luqueno"
$ws.Cells.Item(143, 7).Value = "Nacional Asuncion"
$ws.Cells.Item(143, 8).Value = 1
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = "D"
$ws.Cells.Item(143, 11).Value = 2.75
$ws.Cells.Item(143, 12).Value = 3.2
$ws.Cells.Item(143, 13).Value = 2.4
$ws.Cells.Item(143, 14).Value = 2.75
$ws.Cells.Item(143, 15).Value = 3.1
$ws.Cells.Item(143, 16).Value = 2.45
$ws.Cells.Item(143, 17).Value = 0.25
$ws.Cells.Item(143, 18).Value = 1.75
$ws.Cells.Item(143, 19).Value = 2.05
$ws.Cells.Item(143, 20).Value = 2.25
$ws.Cells.Item(143, 21).Value = 2
$ws.Cells.Item(143, 22).Value = 1.8
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = 2.1
$ws.Cells.Item(143, 25).Value = -1
$ws.Cells.Item(143, 26).Value = 0.375
$ws.Cells.Item(143, 27).Value = -0.5
$ws.Cells.Item(143, 28).Value = -0.5
$ws.Cells.Item(143, 29).Value = 0.4

# Row 145
$ws.Cells.Item(145, 2).Value = 7493311
$ws.Cells.Item(145, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(145, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(145, 5).Value = 45261.8125
$ws.Cells.Item(145, 6).Value = "General Caballero JLM"
$ws.Cells.Item(145, 7).Value = "Olimpia Asuncion"
$ws.Cells.Item(145, 8).Value = 0
$ws.Cells.Item(145, 9).Value = 1
$ws.Cells.Item(145, 10).Value = "A"
$ws.Cells.Item(145, 11).Value = 3.4
$ws.Cells.Item(145, 12).Value = 3.3
$ws.Cells.Item(145, 13).Value = 2
$ws.Cells.Item(145, 14).Value = 3.2
$ws.Cells.Item(145, 15).Value = 3.25
$ws.Cells.Item(145, 16).Value = 2.1
$ws.Cells.Item(145, 17).Value = 0.25
$ws.Cells.Item(145, 18).Value = 1.95
$ws.Cells.Item(145, 19).Value = 1.85
$ws.Cells.Item(145, 20).Value = 2.25
$ws.Cells.Item(145, 21).Value = 1.775
$ws.Cells.Item(145, 22).Value = 2.025
$ws.Cells.Item(145, 23).Value = -1
$ws.Cells.Item(145, 24).Value = -1
$ws.Cells.Item(145, 25).Value = 1.1
$ws.Cells.Item(145, 26).Value = -1
$ws.Cells.Item(145, 27).Value = 0.8500000000000001
$ws.Cells.Item(145, 28).Value = -1
$ws.Cells.Item(145, 29).Value = 1.025

# Row 188
$ws.Cells.Item(188, 2).Value = 7609675
$ws.Cells.Item(188, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(188, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(188, 5).Value = 45353.75
$ws.Cells.Item(188, 6).Value = "Sportivo Ameliano"
$ws.Cells.Item(188, 7).Value = "2 de Mayo"
$ws.Cells.Item(188, 8).Value = 0
$ws.Cells.Item(188, 9).Value = 2
$ws.Cells.Item(188, 10).Value = "A"
$ws.Cells.Item(188, 11).Value = 2.2
$ws.Cells.Item(188, 12).Value = 3.4
$ws.Cells.Item(188, 13).Value = 3
$ws.Cells.Item(188, 14).Value = 2.15
$ws.Cells.Item(188, 15).Value = 3.5
$ws.Cells.Item(188, 16).Value = 3
$ws.Cells.Item(188, 17).Value = -0.25
$ws.Cells.Item(188, 18).Value = 1.875
$ws.Cells.Item(188, 19).Value = 1.925
$ws.Cells.Item(188, 20).Value = 2.5
$ws.Cells.Item(188, 21).Value = 1.875
$ws.Cells.Item(188, 22).Value = 1.925
$ws.Cells.Item(188, 23).Value = -1
$ws.Cells.Item(188, 24).Value = -1
$ws.Cells.Item(188, 25).Value = 2
$ws.Cells.Item(188, 26).Value = -1
$ws.Cells.Item(188, 27).Value = 0.925
$ws.Cells.Item(188, 28).Value = -1
$ws.Cells.Item(188, 29).Value = 0.925

# Row 189
$ws.Cells.Item(189, 2).Value = 7609140
$ws.Cells.Item(189, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(189, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(189, 5).Value = 45354.75
$ws.Cells.Item(189, 6).Value = "General Caballero JLM"
$ws.Cells.Item(189, 7).Value = "Sportivo Trinidense"
$ws.Cells.Item(189, 11).Value = 2.375
$ws.Cells.Item(189, 12).Value = 3.3
$ws.Cells.Item(189, 13).Value = 2.8
$ws.Cells.Item(189, 14).Value = 1.85
$ws.Cells.Item(189, 15).Value = 3.5
$ws.Cells.Item(189, 16).Value = 3.8
$ws.Cells.Item(189, 17).Value = -0.5
$ws.Cells.Item(189, 18).Value = 1.9
$ws.Cells.Item(189, 19).Value = 1.9
$ws.Cells.Item(189, 20).Value = 2.5
$ws.Cells.Item(189, 21).Value = 1.9
$ws.Cells.Item(189, 22).Value = 1.9
$ws.Cells.Item(189, 23).Value = 0
$ws.Cells.Item(189, 24).Value = 0
$ws.Cells.Item(189, 25).Value = 0
$ws.Cells.Item(189, 26).Value = 0
$ws.Cells.Item(189, 27).Value = 0

# Row 190
$ws.Cells.Item(190, 2).Value = 7609193
$ws.Cells.Item(190, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(190, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(190, 5).Value = 45354.85416666666
$ws.Cells.Item(190, 6).Value = "Nacional Asuncion"
$ws.Cells.Item(190, 7).Value = "Libertad Asuncion"
$ws.Cells.Item(190, 11).Value = 4.75
$ws.Cells.Item(190, 12).Value = 3.6
$ws.Cells.Item(190, 13).Value = 1.666
$ws.Cells.Item(190, 14).Value = 4.75
$ws.Cells.Item(190, 15).Value = 3.6
$ws.Cells.Item(190, 16).Value = 1.666
$ws.Cells.Item(190, 17).Value = 0.75
$ws.Cells.Item(190, 18).Value = 1.925
$ws.Cells.Item(190, 19).Value = 1.875
$ws.Cells.Item(190, 20).Value = 2.25
$ws.Cells.Item(190, 21).Value = 1.8
$ws.Cells.Item(190, 22).Value = 2
$ws.Cells.Item(190, 23).Value = 0
$ws.Cells.Item(190, 24).Value = 0
$ws.Cells.Item(190, 25).Value = 0
$ws.Cells.Item(190, 26).Value = 0
$ws.Cells.Item(190, 27).Value = 0

# Row 191
$ws.Cells.Item(191, 2).Value = 7609194
$ws.Cells.Item(191, 3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(191, 4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(191, 5).Value = 45355.8125
$ws.Cells.Item(191, 6).Value = "Tacuary"
$ws.Cells.Item(191, 7).Value = "Sol de America"
$ws.Cells.Item(191, 11).Value = 2
$ws.Cells.Item(191, 12).Value = 3.2
$ws.Cells.Item(191, 13).Value = 3.6
$ws.Cells.Item(191, 14).Value = 2.15
$ws.Cells.Item(191, 15).Value = 3.1
$ws.Cells.Item(191, 16).Value = 3.25
$ws.Cells.Item(191, 17).Value = -0.25
$ws.Cells.Item(191, 18).Value = 1.9
$ws.Cells.Item(191, 19).Value = 1.9
$ws.Cells.Item(191, 20).Value = 2.25
$ws.Cells.Item(191, 21).Value = 1.9
$ws.Cells.Item(191, 22).Value = 1.9
$ws.Cells.Item(191, 23).Value = 0
$ws.Cells.Item(191, 24).Value = 0
$ws.Cells.Item(191, 25).Value = 0
$ws.Cells.Item(191, 26).Value = 0
$ws.Cells.Item(191, 27).Value = 0

# Row 192 content is obsolete after the shift; delete the row entirely
$ws.Rows.Item(192).Delete()
